# edit.ps1 - apply the "fodo.pptx" revision:
#   1. bump the cached datetimeFigureOut placeholder text (slide master + all
#      slide layouts) from 5/1/2015 -> 5/2/2015
#   2. on slide 1, shrink the "Initial beam parameters" textbox and tweak its
#      wording (merge "2 " + "micron" runs, change Twiss beta value, reword
#      the "verify that" line)
#   3. on slide 1, merge split runs in the "Objective 4" bullet and the
#      "Hint:" line

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder: slide master + every slide layout.
#    NB: shapes are addressed with explicit loops (no helper function) --
#    passing the Shapes COM collection into a PowerShell function here loses
#    the collection (Count comes back 0), so everything is inlined.
# ---------------------------------------------------------------------------
$newDate = "5/2/2015"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 - "TextBox 6" (the "Initial beam parameters:" box)
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$paramsBox = $slide1.Shapes.Item("TextBox 6")

# narrower box (height is unchanged)
$paramsBox.Width = 449.7192913385827

$tr = $paramsBox.TextFrame.TextRange

# "2 " + "micron" -> single run "2 micron"
$emittanceLine = $tr.Lines(3, 1)
$emittanceLine.Characters(33, 8).Text = "2 micron"

# "0.02 " -> "2.5 " (keeps the trailing "m" as its own run)
$twissLine = $tr.Lines(4, 1)
$twissLine.Characters(25, 5).Text = "2.5 "

# "verify that " -> "verify:  "
$verifyLine = $tr.Lines(6, 1)
$verifyLine.Characters(1, 12).Text = "verify:  "

# ---------------------------------------------------------------------------
# 3) Slide 1 - "Objective 4" bullet: merge ":   " + "Create a new beamline..."
# ---------------------------------------------------------------------------
$obj4Box = $slide1.Shapes.Item(8)
$obj4Tr = $obj4Box.TextFrame.TextRange
$obj4Tr.Characters(13, 55).Text = ":   Create a new beamline, consisting of 10 FODO cells "

# ---------------------------------------------------------------------------
# 4) Slide 1 - "TextBox 11" (Hint): merge ": " + " " -> ":  "
# ---------------------------------------------------------------------------
$hintBox = $slide1.Shapes.Item("TextBox 11")
$hintLine = $hintBox.TextFrame.TextRange.Lines(3, 1)
$hintLine.Characters(5, 3).Text = ":  "

Write-Host "fodo.pptx edits applied"
